$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock quantity values
$ws.Range("C2").Value = 4678
$ws.Range("C4").Value = 699

# Update the active selection on the sheet
$ws.Range("D2").Select()
